$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text so numeric-looking strings (e.g. "1.00",
# "0.515") keep their original text representation instead of being coerced
# into numbers (which would drop trailing zeros / change formatting).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows with refreshed market data
Set-TextValue $ws.Range("D2") "29.612.27"
Set-TextValue $ws.Range("E2") "  +1.84%  "
Set-TextValue $ws.Range("D3") "1.598.48"
Set-TextValue $ws.Range("E3") "  +1.13%  "
Set-TextValue $ws.Range("E4") "  +0.42%  "
Set-TextValue $ws.Range("D5") "212.19"
Set-TextValue $ws.Range("E5") "  +0.18%  "
Set-TextValue $ws.Range("D6") "0.515"
Set-TextValue $ws.Range("E6") "  -0.83%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.45%  "
Set-TextValue $ws.Range("D8") "26.82"
Set-TextValue $ws.Range("E8") "  +4.75%  "
Set-TextValue $ws.Range("D9") "43.81"
Set-TextValue $ws.Range("E9") "  -1.65%  "
Set-TextValue $ws.Range("E10") "  +1.23%  "
Set-TextValue $ws.Range("E11") "  +0.90%  "
Set-TextValue $ws.Range("D12") "0.0909"
Set-TextValue $ws.Range("E12") "  +1.08%  "
Set-TextValue $ws.Range("D13") "1.826.98"
Set-TextValue $ws.Range("E13") "  +1.08%  "
Set-TextValue $ws.Range("D14") "1.601.13"
Set-TextValue $ws.Range("E14") "  +1.10%  "
Set-TextValue $ws.Range("D15") "29.615.08"
Set-TextValue $ws.Range("E15") "  +1.76%  "
Set-TextValue $ws.Range("E16") "  +3.08%  "
Set-TextValue $ws.Range("E17") "  +0.87%  "
Set-TextValue $ws.Range("D18") "63.89"
Set-TextValue $ws.Range("E18") "  +2.36%  "
Set-TextValue $ws.Range("D19") "241.63"
Set-TextValue $ws.Range("E19") "  +1.91%  "
Set-TextValue $ws.Range("D20") "7.61"
Set-TextValue $ws.Range("E20") "  +2.65%  "
Set-TextValue $ws.Range("E21") "  +0.23%  "
Set-TextValue $ws.Range("E22") "  +0.38%  "
Set-TextValue $ws.Range("E23") "  -0.40%  "
Set-TextValue $ws.Range("E24") "  +0.39%  "
Set-TextValue $ws.Range("E25") "  -0.67%  "
Set-TextValue $ws.Range("D26") "154.58"
Set-TextValue $ws.Range("E26") "  +0.95%  "
Set-TextValue $ws.Range("E27") "  +1.97%  "
Set-TextValue $ws.Range("D28") "0.109"
Set-TextValue $ws.Range("E28") "  +0.00%  "
Set-TextValue $ws.Range("E29") "  +1.23%  "
Set-TextValue $ws.Range("E30") "  +0.38%  "
Set-TextValue $ws.Range("D31") "0.0476"
Set-TextValue $ws.Range("E31") "  +2.68%  "
Set-TextValue $ws.Range("D32") "1.06"
Set-TextValue $ws.Range("E32") "  +0.62%  "
Set-TextValue $ws.Range("E33") "  +0.31%  "
Set-TextValue $ws.Range("E34") "  +3.19%  "
Set-TextValue $ws.Range("D35") "1.430.16"
Set-TextValue $ws.Range("E35") "  +0.79%  "
Set-TextValue $ws.Range("E36") "  +2.12%  "
Set-TextValue $ws.Range("E37") "  -1.60%  "
Set-TextValue $ws.Range("D38") "2.87"
Set-TextValue $ws.Range("E38") "  +3.63%  "
Set-TextValue $ws.Range("E39") "  +0.25%  "
Set-TextValue $ws.Range("E40") "  +1.64%  "
Set-TextValue $ws.Range("D41") "0.541"
Set-TextValue $ws.Range("E41") "  +3.07%  "
Set-TextValue $ws.Range("D44") "0.0494"
Set-TextValue $ws.Range("E44") "  +6.81%  "
Set-TextValue $ws.Range("D45") "0.802"
Set-TextValue $ws.Range("E45") "  +2.26%  "
Set-TextValue $ws.Range("D46") "1.00"
Set-TextValue $ws.Range("E46") "  +0.39%  "
Set-TextValue $ws.Range("D47") "0.984"
Set-TextValue $ws.Range("E47") "  +15.33%  "
Set-TextValue $ws.Range("D48") "65.72"
Set-TextValue $ws.Range("E48") "  +1.74%  "
Set-TextValue $ws.Range("D49") "5.32"
Set-TextValue $ws.Range("E49") "  -0.12%  "
Set-TextValue $ws.Range("D50") "1.738.52"
Set-TextValue $ws.Range("E50") "  +1.10%  "
Set-TextValue $ws.Range("D51") "86.09"
Set-TextValue $ws.Range("E51") "  +0.52%  "

# Rows 42 and 43 swap ranking order (BitcoinSV moves up to 40, RenderToken drops to 41)
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D42") "54.52"
Set-TextValue $ws.Range("E42") "  +3.31%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "1.96"
Set-TextValue $ws.Range("E43") "  +0.94%  "
